# Adds simulated patient data for three new patients/branches (1568, 1569,
# 1570) to the "data" worksheet, continuing directly after the existing
# last row (493).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Shared-string text values used throughout the table
# (mirrors xl/sharedStrings.xml ordering)
$PER   = "PER"
$VIS   = "VIS"
$D_ERA = "D_ERA"
$D_EXP = "D_EXP"
$Daily = "Daily"
$Twice = "2 times daily"
$Four  = "4 times daily"
$Null_ = "null"

# Each entry describes one new row with its column values
# (R=row, then column letter => value). Missing keys mean the cell is left blank.
$rows = @(
    @{ R=494; A=$PER;   B=1568; C=1960;                                      K=1;  L=1 },
    @{ R=495; A=$D_ERA; B=1568; D=1398937;   E=39491; F=39565 },
    @{ R=496; A=$D_ERA; B=1568; D=902427;    E=39491; F=39565 },
    @{ R=497; A=$VIS;   B=1568;              E=39491; F=39565; H=9201 },
    @{ R=498; A=$D_EXP; B=1568; D=40223504;  E=39491; F=39565; I=30; J=30; M=$Twice; N=$Null_; O=$Null_ },
    @{ R=499; A=$D_EXP; B=1568; D=1594707;   E=39491; F=39565; I=30; J=30; M=$Twice; N=$Null_; O=$Null_ },

    @{ R=500; A=$PER;   B=1569; C=1966;                                      K=1;  L=1 },
    @{ R=501; A=$D_ERA; B=1569; D=1398937;   E=39491; F=39565 },
    @{ R=502; A=$D_ERA; B=1569; D=902427;    E=39491; F=39565 },
    @{ R=503; A=$VIS;   B=1569;              E=39491; F=39565; H=9201 },
    @{ R=504; A=$D_EXP; B=1569; D=19079775;  E=39491; F=39565; I=20; J=30; M=$Four; N=$Null_; O=$Null_ },
    @{ R=505; A=$D_EXP; B=1569; D=902489;    E=39491; F=39565; I=5;  J=30; M=$Four; N=$Null_; O=$Null_ },

    @{ R=506; A=$PER;   B=1570; C=1950;                                      K=1;  L=1 },
    @{ R=507; A=$D_ERA; B=1570; D=1398937;   E=39491; F=39565 },
    @{ R=508; A=$D_ERA; B=1570; D=950370;    E=39491; F=39565 },
    @{ R=509; A=$VIS;   B=1570;              E=39491; F=39565; H=9201 },
    @{ R=510; A=$D_EXP; B=1570; D=40223506;  E=39491; F=39565; I=10; J=10; M=$Daily; N=$Null_; O=$Null_ },
    @{ R=511; A=$D_EXP; B=1570; D=43219718;  E=39491; F=39565; I=30; J=30; M=$Twice; N=$Null_; O=$Null_ }
)

# Use an existing formatted date cell as the format source so pasted date
# cells reuse the workbook's existing "short date" style instead of Excel
# registering a brand-new (duplicate) style for every cell.
$dateFormatSrc = $ws.Cells.Item(3, 5)

foreach ($row in $rows) {
    $r = $row.R

    $ws.Cells.Item($r, 1).Value = $row.A

    if ($row.ContainsKey("B")) { $ws.Cells.Item($r, 2).Value = $row.B }
    if ($row.ContainsKey("C")) { $ws.Cells.Item($r, 3).Value = $row.C }
    if ($row.ContainsKey("D")) { $ws.Cells.Item($r, 4).Value = $row.D }

    if ($row.ContainsKey("E")) {
        $cell = $ws.Cells.Item($r, 5)
        $cell.Value = $row.E
        $dateFormatSrc.Copy()
        $cell.PasteSpecial(-4122)
    }
    if ($row.ContainsKey("F")) {
        $cell = $ws.Cells.Item($r, 6)
        $cell.Value = $row.F
        $dateFormatSrc.Copy()
        $cell.PasteSpecial(-4122)
    }

    if ($row.ContainsKey("H")) { $ws.Cells.Item($r, 8).Value = $row.H }
    if ($row.ContainsKey("I")) { $ws.Cells.Item($r, 9).Value = $row.I }
    if ($row.ContainsKey("J")) { $ws.Cells.Item($r, 10).Value = $row.J }
    if ($row.ContainsKey("K")) { $ws.Cells.Item($r, 11).Value = $row.K }
    if ($row.ContainsKey("L")) { $ws.Cells.Item($r, 12).Value = $row.L }
    if ($row.ContainsKey("M")) { $ws.Cells.Item($r, 13).Value = $row.M }
    if ($row.ContainsKey("N")) { $ws.Cells.Item($r, 14).Value = $row.N }
    if ($row.ContainsKey("O")) { $ws.Cells.Item($r, 15).Value = $row.O }
}

# Update the selection to reflect scrolling down to the newly added rows,
# same as what the author ended up with after typing the data (the sheet
# already has the header row frozen, which is left untouched).
$ws.Activate() | Out-Null
$ws.Range("L507").Select() | Out-Null
